$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells keep their original text formatting so
# values such as "28.407.05" or "0.00001135" are not reinterpreted as numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '28.407.05'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '1.820.69'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '316.25'
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').Value = '0.5096'
$ws.Range('E7').Value = '  -4.17%  '
$ws.Range('D8').Value = '0.3922'
$ws.Range('E8').Value = '  -3.30%  '
$ws.Range('D9').Value = '0.07760'
$ws.Range('E9').Value = '  +1.75%  '
$ws.Range('D10').Value = '41.78'
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('D11').Value = '1.109'
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').Value = '20.90'
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').Value = '1.001'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('E14').Value = '  -1.84%  '
$ws.Range('D15').Value = '7.468'
$ws.Range('E15').Value = '  -1.96%  '
$ws.Range('D16').Value = '1.810.85'
$ws.Range('E16').Value = '  -0.92%  '
$ws.Range('D17').Value = '0.00001135'
$ws.Range('E17').Value = '  +5.60%  '
$ws.Range('D18').Value = '92.42'
$ws.Range('E18').Value = '  +3.35%  '
$ws.Range('D19').Value = '0.06619'
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').Value = '17.71'
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').Value = '6.082'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').Value = '28.442.01'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = '11.27'
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').Value = '2.245'
$ws.Range('E25').Value = '  +4.54%  '
$ws.Range('D26').Value = '21.07'
$ws.Range('E26').Value = '  +2.40%  '
$ws.Range('D27').Value = '2.029.74'
$ws.Range('E27').Value = '  -0.49%  '
$ws.Range('D28').Value = '155.27'
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('D29').Value = '2.408'
$ws.Range('E29').Value = '  -2.96%  '
$ws.Range('D30').Value = '125.18'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('D31').Value = '0.1100'
$ws.Range('E31').Value = '  +0.69%  '
$ws.Range('D32').Value = '1.101'
$ws.Range('E32').Value = '  -2.33%  '
$ws.Range('E33').Value = '  -0.52%  '
$ws.Range('D34').Value = '3.645'
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('D35').Value = '0.07035'
$ws.Range('E35').Value = '  -1.59%  '
$ws.Range('D36').Value = '0.2207'
$ws.Range('E36').Value = '  -2.46%  '
$ws.Range('D37').Value = '0.02321'
$ws.Range('E37').Value = '  -0.94%  '
$ws.Range('D38').Value = '5.191'
$ws.Range('E38').Value = '  -0.52%  '
$ws.Range('D39').Value = '8.748'
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('D40').Value = '0.6257'
$ws.Range('E40').Value = '  -0.33%  '
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('D42').Value = '1.174'
$ws.Range('E42').Value = '  -1.11%  '
$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E44').Value = '  -1.08%  '
$ws.Range('D45').Value = '13.37'
$ws.Range('E45').Value = '  -0.92%  '
$ws.Range('D46').Value = '3.730'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('D47').Value = '0.5873'
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('D48').Value = '124.20'
$ws.Range('E48').Value = '  -1.42%  '
$ws.Range('D49').Value = '1.974'
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('D51').Value = '0.06896'
$ws.Range('E51').Value = '  -0.05%  '
